$d = $word.ActiveDocument

function New-FlatOpcXml($bodyInnerXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyInnerXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

$rFonts = '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>'

# ---------------------------------------------------------------------------
# Step 1 (applied first, works on the *last* paragraph so earlier offsets in
# the document remain valid): replace the "users" list paragraph (which also
# drops its spellStart/spellEnd proofErr wrapper) and append all of the new
# paragraphs describing the user / moderator rights after it.
# ---------------------------------------------------------------------------
$usersPara = $d.Paragraphs(14)
if ($usersPara.Range.Text.TrimEnd([char]13) -ne "users") {
    throw "Unexpected paragraph 14 content: " + $usersPara.Range.Text
}

$tailXml = @"
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="13"/>
    </w:numPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>u</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>sers</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t xml:space="preserve">Pour </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>le user</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t xml:space="preserve"> u</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>n moyen de :</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="13"/>
    </w:numPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>S</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>'inscrire</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="13"/>
    </w:numPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>Se connecter</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="13"/>
    </w:numPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>D'accéder à un espace user</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="13"/>
    </w:numPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>D'ajouter des commentaires</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t xml:space="preserve">Pour le modérateur s'ajoute : </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="13"/>
    </w:numPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>Le droit d'ajouter des éventments</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
</w:p>
"@

$usersPara.Range.InsertXML((New-FlatOpcXml $tailXml)) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: drop the _GoBack bookmark that used to sit in the "Des tables : "
# paragraph.
# ---------------------------------------------------------------------------
$tablesPara = $d.Paragraphs(9)
if ($tablesPara.Range.Text.TrimEnd([char]13) -ne "Des tables : ") {
    throw "Unexpected paragraph 9 content: " + $tablesPara.Range.Text
}

$tablesXml = @"
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t xml:space="preserve">Des tables : </w:t>
  </w:r>
</w:p>
"@

$tablesPara.Range.InsertXML((New-FlatOpcXml $tablesXml)) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: insert a new paragraph "Des formulaires pour communiquer avec la
# base de données" right before "Une base de données (lesevenements)".
# ---------------------------------------------------------------------------
$dbPara = $d.Paragraphs(8)
if ($dbPara.Range.Text.TrimEnd([char]13) -ne "Une base de données (lesevenements)") {
    throw "Unexpected paragraph 8 content: " + $dbPara.Range.Text
}

$dbXml = @"
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>Des formulaires pour communiquer avec la base de données</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      $rFonts
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      $rFonts
    </w:rPr>
    <w:t>Une base de données (lesevenements)</w:t>
  </w:r>
</w:p>
"@

$dbPara.Range.InsertXML((New-FlatOpcXml $dbXml)) | Out-Null

Write-Output ("Done. Paragraph count = " + $d.Paragraphs.Count)
